# Refactor code scripts added.
# Fill in the Transaction Reference / Dr Account / Status results that were
# recorded for the latest regression run (rows 3-16 of Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3 -----------------------------------------------------------
$ws.Range("E3").Value = "081023113057825P"
$ws.Range("M3").Value = "PASS"

# --- Row 4 -----------------------------------------------------------
$ws.Range("E4").Value = "081023113314459P"
$ws.Range("M4").Value = "PASS"

# --- Row 5 -----------------------------------------------------------
$ws.Range("E5").Value = "081023120223568P"
$ws.Range("M5").Value = "PASS"

# --- Row 9 -----------------------------------------------------------
$ws.Range("E9").Value = "081023153817961P"
$ws.Range("M9").Value = "PASS"

# --- Row 10 ------------------------------------------------------------
$ws.Range("E10").Value = "081023154033571P"
# F10 needs the same Dr Account number already used elsewhere (shared
# string "1501200131929001"). Assigning the plain numeric-looking text
# directly would make Excel store it as a real number (and would stamp a
# new quote-prefixed style onto the cell), so instead copy the value from
# a cell that already holds it as text and paste values-only, which keeps
# F10's own style (s="42") untouched.
$ws.Range("F9").Copy() | Out-Null
$ws.Range("F10").PasteSpecial(-4163) | Out-Null

# --- Row 11 ------------------------------------------------------------
$ws.Range("E11").Value = "081023154318114P"
$ws.Range("M11").Value = "PASS"

# --- Row 15 ------------------------------------------------------------
$ws.Range("E15").Value = "041023104739348P"

# --- Row 16 ------------------------------------------------------------
$ws.Range("E16").Value = "091023114001997P"
$ws.Range("F9").Copy() | Out-Null
$ws.Range("F16").PasteSpecial(-4163) | Out-Null
$ws.Range("M16").Value = "PASS"

$excel.CutCopyMode = 0
